$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NIK typo fix, applies to every row that shared the old value
$ws.Range("C2:C7").Value = "EN-4-046"

# Row 2 (Bulan = 6 / June) is for a different employee
$ws.Range("D2").Value = "Agus Priyanto"

# Terlambat (late) count for row 2 updated
$ws.Range("G2").Value = 2

# Selection moved to G3
$ws.Range("G3").Select()
